# Add Wins / Losses / Ties columns (AD, AE, AF) with the team's season record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): add new column headers
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting of the other header cells (bold, bordered, centered style)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Determine the last used data row (header is row 1, data starts at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Fill in the season record for every player row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 91
    $ws.Cells.Item($r, 31).Value = 71
    $ws.Cells.Item($r, 32).Value = 0
}
